# Generate Report for Handoff
#
# Adds two new "Ready for handoff" rows (for newly-handed-off files
# c8da94ef-7763-4898-95c5-6f57c0071fd8 and d9d5826d-78db-42bb-a25c-391c27bc5a40)
# to every worksheet in the localization-status report:
#   - Overview : row 6 / row 7
#   - zh-cn    : row 6 / row 7
#   - de-de    : row 6 / row 7

$wb = $excel.ActiveWorkbook

# Hyperlink colour used throughout the workbook for the "HyperLink" cell style
# (font color FF6495ED, underlined) - reapplied to every new hyperlink cell so
# the look matches the existing rows.
$hyperlinkColor = 15570276   # BGR value of RGB(100,149,237) / #6495ED

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

function Add-Link($ws, $cellRange, $text, $url) {
    $cellRange.Value = $text
    $ws.Hyperlinks.Add($cellRange, $url, [Type]::Missing, [Type]::Missing, $text) | Out-Null
    Style-AsHyperlink $cellRange
}

# ---------------------------------------------------------------------------
# New file identifiers
# ---------------------------------------------------------------------------
$guid1 = "c8da94ef-7763-4898-95c5-6f57c0071fd8"
$guid2 = "d9d5826d-78db-42bb-a25c-391c27bc5a40"

$sha1 = "196dcaf38e049fbb599ac88b9ecdbe0f263b8f2d"
$sha2 = "7a8ca4b25bc6ffaaf81728f3bf2a2213289309df"

$status = "Ready for handoff"

$dateOverview1 = "2016-33-21 00:33:35"
$dateOverview2 = "2016-33-21 00:33:35"

$dateZh = "2016-03-21 00:33:31"
$dateDe = "2016-03-21 00:33:35"

$emptyDate = "0001-01-01 00:00:00"
$reason = "Include"

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e"

# Pre-built text / URL strings (built ahead of time so that no call site has
# to chain two separate "-f" expressions as adjacent arguments).
$md1Text = $guid1 + ".md"
$md1Url = $baseUrl + "/" + $guid1 + ".md"

$md2Text = $guid2 + ".md"
$md2Url = $baseUrl + "/" + $guid2 + ".md"

$zh1Text = $guid1 + "." + $sha1 + ".zh-cn.xlf"
$zh1Url = $baseUrl + "/" + $zh1Text

$zh2Text = $guid2 + "." + $sha2 + ".zh-cn.xlf"
$zh2Url = $baseUrl + "/" + $zh2Text

$de1Text = $guid1 + "." + $sha1 + ".de-de.xlf"
$de1Url = $baseUrl + "/" + $de1Text

$de2Text = $guid2 + "." + $sha2 + ".de-de.xlf"
$de2Url = $baseUrl + "/" + $de2Text

$dotMdText = ".md"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# --- Row 6 : guid1 -----------------------------------------------------
Add-Link $wsOverview $wsOverview.Range("A6") $md1Text $md1Url
$wsOverview.Range("B6").Value = $status
$wsOverview.Range("C6").Value = $status
$wsOverview.Range("D6").Value = $dateOverview1

# --- Row 7 : guid2 -----------------------------------------------------
Add-Link $wsOverview $wsOverview.Range("A7") $md2Text $md2Url
$wsOverview.Range("B7").Value = $status
$wsOverview.Range("C7").Value = $status
$wsOverview.Range("D7").Value = $dateOverview2

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# --- Row 6 : guid1 -------------------------------------------------------
Add-Link $wsZh $wsZh.Range("A6") $md1Text $md1Url
Add-Link $wsZh $wsZh.Range("B6") $dotMdText $md1Url
$wsZh.Range("C6").Value = $status
Add-Link $wsZh $wsZh.Range("D6") $zh1Text $zh1Url
$wsZh.Range("E6").Value = $dateZh
$wsZh.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H6").Value = $emptyDate
$wsZh.Range("I6").Value = $reason

# --- Row 7 : guid2 -------------------------------------------------------
Add-Link $wsZh $wsZh.Range("A7") $md2Text $md2Url
Add-Link $wsZh $wsZh.Range("B7") $dotMdText $md2Url
$wsZh.Range("C7").Value = $status
Add-Link $wsZh $wsZh.Range("D7") $zh2Text $zh2Url
$wsZh.Range("E7").Value = $dateZh
$wsZh.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H7").Value = $emptyDate
$wsZh.Range("I7").Value = $reason

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

# --- Row 6 : guid1 -------------------------------------------------------
Add-Link $wsDe $wsDe.Range("A6") $md1Text $md1Url
Add-Link $wsDe $wsDe.Range("B6") $dotMdText $md1Url
$wsDe.Range("C6").Value = $status
Add-Link $wsDe $wsDe.Range("D6") $de1Text $de1Url
$wsDe.Range("E6").Value = $dateDe
$wsDe.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H6").Value = $emptyDate
$wsDe.Range("I6").Value = $reason

# --- Row 7 : guid2 -------------------------------------------------------
Add-Link $wsDe $wsDe.Range("A7") $md2Text $md2Url
Add-Link $wsDe $wsDe.Range("B7") $dotMdText $md2Url
$wsDe.Range("C7").Value = $status
Add-Link $wsDe $wsDe.Range("D7") $de2Text $de2Url
$wsDe.Range("E7").Value = $dateDe
$wsDe.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H7").Value = $emptyDate
$wsDe.Range("I7").Value = $reason
